$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1933701657458564
$ws.Range("C2").Value = 0.569060773480663
$ws.Range("J2").Value = 0.01657458563535912
$ws.Range("P2").Value = 0.143646408839779
$ws.Range("S2").Value = 0.07734806629834254
$ws.Range("C3").Value = 0.03703703703703703
$ws.Range("J3").Value = 0.03703703703703703
$ws.Range("P3").Value = 0.6574074074074074
$ws.Range("S3").Value = 0.2685185185185185
$ws.Range("J4").Value = 0.08333333333333333
$ws.Range("P4").Value = 0.7083333333333334
$ws.Range("S4").Value = 0.2083333333333333
$ws.Range("B6").Value = 0.09142857142857143
$ws.Range("D6").Value = 0.01714285714285714
$ws.Range("F6").Value = 0.04571428571428571
$ws.Range("J6").Value = 0.2342857142857143
$ws.Range("O6").Value = 0.03428571428571429
$ws.Range("Q6").Value = 0.1371428571428571
$ws.Range("R6").Value = 0.07428571428571429
$ws.Range("S6").Value = 0.3657142857142857
$ws.Range("B7").Value = 0.1043478260869565
$ws.Range("D7").Value = 0.01739130434782609
$ws.Range("F7").Value = 0.06956521739130435
$ws.Range("J7").Value = 0.1217391304347826
$ws.Range("O7").Value = 0.01739130434782609
$ws.Range("Q7").Value = 0.1739130434782609
$ws.Range("R7").Value = 0.1130434782608696
$ws.Range("S7").Value = 0.3826086956521739
$ws.Range("B8").Value = 0.06049822064056939
$ws.Range("D8").Value = 0.01779359430604982
$ws.Range("F8").Value = 0.05338078291814947
$ws.Range("J8").Value = 0.103202846975089
$ws.Range("O8").Value = 0.01423487544483986
$ws.Range("Q8").Value = 0.2241992882562278
$ws.Range("R8").Value = 0.103202846975089
$ws.Range("S8").Value = 0.4234875444839858
$ws.Range("B9").Value = 0.06622516556291391
$ws.Range("D9").Value = 0.006622516556291391
$ws.Range("F9").Value = 0.06622516556291391
$ws.Range("J9").Value = 0.07947019867549669
$ws.Range("O9").Value = 0.03973509933774835
$ws.Range("Q9").Value = 0.1456953642384106
$ws.Range("R9").Value = 0.1324503311258278
$ws.Range("S9").Value = 0.4635761589403973
$ws.Range("B10").Value = 0.08670520231213873
$ws.Range("D10").Value = 0.01541425818882466
$ws.Range("E10").Value = 0.0009633911368015414
$ws.Range("F10").Value = 0.07418111753371869
$ws.Range("J10").Value = 0.1021194605009634
$ws.Range("O10").Value = 0.02119460500963391
$ws.Range("Q10").Value = 0.2235067437379576
$ws.Range("R10").Value = 0.1107899807321773
$ws.Range("S10").Value = 0.3651252408477842
$ws.Range("G11").Value = 0.1141304347826087
$ws.Range("J11").Value = 0.108695652173913
$ws.Range("K11").Value = 0.1739130434782609
$ws.Range("L11").Value = 0.5543478260869565
$ws.Range("S11").Value = 0.04891304347826087
$ws.Range("G12").Value = 0.7592592592592593
$ws.Range("J12").Value = 0.1481481481481481
$ws.Range("K12").Value = 0.009259259259259259
$ws.Range("L12").Value = 0.03703703703703703
$ws.Range("S12").Value = 0.04629629629629629
$ws.Range("G13").Value = 0.5185185185185185
$ws.Range("J13").Value = 0.3703703703703703
$ws.Range("S13").Value = 0.1111111111111111
$ws.Range("F15").Value = 0.04210526315789474
$ws.Range("H15").Value = 0.1631578947368421
$ws.Range("I15").Value = 0.08421052631578947
$ws.Range("J15").Value = 0.3736842105263158
$ws.Range("K15").Value = 0.03684210526315789
$ws.Range("M15").Value = 0.005263157894736842
$ws.Range("O15").Value = 0.05263157894736842
$ws.Range("S15").Value = 0.2421052631578947
$ws.Range("F16").Value = 0.01818181818181818
$ws.Range("H16").Value = 0.05454545454545454
$ws.Range("I16").Value = 0.04545454545454546
$ws.Range("J16").Value = 0.6090909090909091
$ws.Range("K16").Value = 0.08181818181818182
$ws.Range("M16").Value = 0.03636363636363636
$ws.Range("O16").Value = 0.09090909090909091
$ws.Range("S16").Value = 0.06363636363636363
$ws.Range("F17").Value = 0.01955307262569832
$ws.Range("H17").Value = 0.1396648044692737
$ws.Range("I17").Value = 0.09217877094972067
$ws.Range("J17").Value = 0.4720670391061452
$ws.Range("K17").Value = 0.09217877094972067
$ws.Range("M17").Value = 0.0223463687150838
$ws.Range("N17").Value = 0.002793296089385475
$ws.Range("O17").Value = 0.06424581005586592
$ws.Range("S17").Value = 0.09497206703910614
$ws.Range("F18").Value = 0.005319148936170213
$ws.Range("H18").Value = 0.1170212765957447
$ws.Range("I18").Value = 0.09574468085106383
$ws.Range("J18").Value = 0.5053191489361702
$ws.Range("K18").Value = 0.0797872340425532
$ws.Range("M18").Value = 0.02659574468085106
$ws.Range("O18").Value = 0.0797872340425532
$ws.Range("S18").Value = 0.09042553191489362
$ws.Range("F19").Value = 0.01541850220264317
$ws.Range("H19").Value = 0.1949339207048458
$ws.Range("I19").Value = 0.08700440528634361
$ws.Range("J19").Value = 0.4295154185022027
$ws.Range("K19").Value = 0.09361233480176212
$ws.Range("M19").Value = 0.013215859030837
$ws.Range("N19").Value = 0.001101321585903084
$ws.Range("O19").Value = 0.08149779735682819
$ws.Range("S19").Value = 0.08370044052863436
